$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "9+14="
$t.Cell(1,2).Range.Text = "37+53="
$t.Cell(1,3).Range.Text = "93-27="
$t.Cell(1,4).Range.Text = "25+38="
$t.Cell(1,5).Range.Text = "46-7="
$t.Cell(2,1).Range.Text = "20+17="
$t.Cell(2,2).Range.Text = "25-8="
$t.Cell(2,3).Range.Text = "2+67="
$t.Cell(2,4).Range.Text = "66+10="
$t.Cell(2,5).Range.Text = "97-65="
$t.Cell(3,1).Range.Text = "75-14="
$t.Cell(3,2).Range.Text = "14+44="
$t.Cell(3,3).Range.Text = "72-56="
$t.Cell(3,4).Range.Text = "18+33="
$t.Cell(3,5).Range.Text = "1+29="
$t.Cell(4,1).Range.Text = "59+27="
$t.Cell(4,2).Range.Text = "74-35="
$t.Cell(4,3).Range.Text = "19-10="
$t.Cell(4,4).Range.Text = "40+54="
$t.Cell(4,5).Range.Text = "43+8="
$t.Cell(5,1).Range.Text = "57-41="
$t.Cell(5,2).Range.Text = "52-26="
$t.Cell(5,3).Range.Text = "68-51="
$t.Cell(5,4).Range.Text = "97-89="
$t.Cell(5,5).Range.Text = "97-12="
$t.Cell(6,1).Range.Text = "71+4="
$t.Cell(6,2).Range.Text = "15+81="
$t.Cell(6,3).Range.Text = "39+23="
$t.Cell(6,4).Range.Text = "42+36="
$t.Cell(6,5).Range.Text = "49-42="
$t.Cell(7,1).Range.Text = "56+4="
$t.Cell(7,2).Range.Text = "86-67="
$t.Cell(7,3).Range.Text = "37-32="
$t.Cell(7,4).Range.Text = "23+39="
$t.Cell(7,5).Range.Text = "57-10="
$t.Cell(8,1).Range.Text = "9+87="
$t.Cell(8,2).Range.Text = "45+37="
$t.Cell(8,3).Range.Text = "86-16="
$t.Cell(8,4).Range.Text = "93-21="
$t.Cell(8,5).Range.Text = "85-65="
$t.Cell(9,1).Range.Text = "75+19="
$t.Cell(9,2).Range.Text = "65+24="
$t.Cell(9,3).Range.Text = "76-5="
$t.Cell(9,4).Range.Text = "79-53="
$t.Cell(9,5).Range.Text = "8+16="
$t.Cell(10,1).Range.Text = "72-64="
$t.Cell(10,2).Range.Text = "68-50="
$t.Cell(10,3).Range.Text = "39+46="
$t.Cell(10,4).Range.Text = "93-52="
$t.Cell(10,5).Range.Text = "15-14="
$t.Cell(11,1).Range.Text = "43-40="
$t.Cell(11,2).Range.Text = "57-46="
$t.Cell(11,3).Range.Text = "6+41="
$t.Cell(11,4).Range.Text = "94-71="
$t.Cell(11,5).Range.Text = "60+37="
$t.Cell(12,1).Range.Text = "21+71="
$t.Cell(12,2).Range.Text = "75-2="
$t.Cell(12,3).Range.Text = "91-76="
$t.Cell(12,4).Range.Text = "47+19="
$t.Cell(12,5).Range.Text = "29+56="
$t.Cell(13,1).Range.Text = "3+91="
$t.Cell(13,2).Range.Text = "58+24="
$t.Cell(13,3).Range.Text = "97-64="
$t.Cell(13,4).Range.Text = "6-3="
$t.Cell(13,5).Range.Text = "36-33="
$t.Cell(14,1).Range.Text = "13+1="
$t.Cell(14,2).Range.Text = "59-13="
$t.Cell(14,3).Range.Text = "88-48="
$t.Cell(14,4).Range.Text = "7+41="
$t.Cell(14,5).Range.Text = "34+31="
$t.Cell(15,1).Range.Text = "74-26="
$t.Cell(15,2).Range.Text = "99+0="
$t.Cell(15,3).Range.Text = "51+7="
$t.Cell(15,4).Range.Text = "44+31="
$t.Cell(15,5).Range.Text = "76-0="
$t.Cell(16,1).Range.Text = "6+33="
$t.Cell(16,2).Range.Text = "99-22="
$t.Cell(16,3).Range.Text = "57-51="
$t.Cell(16,4).Range.Text = "61-39="
$t.Cell(16,5).Range.Text = "31-25="
$t.Cell(17,1).Range.Text = "58-34="
$t.Cell(17,2).Range.Text = "3+21="
$t.Cell(17,3).Range.Text = "3+26="
$t.Cell(17,4).Range.Text = "13+61="
$t.Cell(17,5).Range.Text = "33+25="
$t.Cell(18,1).Range.Text = "61-48="
$t.Cell(18,2).Range.Text = "43-11="
$t.Cell(18,3).Range.Text = "8+81="
$t.Cell(18,4).Range.Text = "53+11="
$t.Cell(18,5).Range.Text = "81-5="
$t.Cell(19,1).Range.Text = "97-52="
$t.Cell(19,2).Range.Text = "42-41="
$t.Cell(19,3).Range.Text = "32+52="
$t.Cell(19,4).Range.Text = "78-2="
$t.Cell(19,5).Range.Text = "38-18="
$t.Cell(20,1).Range.Text = "7+8="
$t.Cell(20,2).Range.Text = "89-46="
$t.Cell(20,3).Range.Text = "19+13="
$t.Cell(20,4).Range.Text = "35-13="
$t.Cell(20,5).Range.Text = "50+35="
